# Fill empty ER sheet ("GEO_RNASEQ") into the 4COM01_RNASeq template.
#
# This adds a new worksheet named "GEO_RNASEQ" at the end of the workbook,
# fills it with the external-reference ("ER") description table that
# mirrors the columns of the main annotation table (sheet "4COM01_RNASeq"),
# sizes its columns, and makes it the active/selected sheet - matching the
# upstream "Fill empty ER sheets into every template (except Imaging)"
# commit.

$wb = $excel.ActiveWorkbook

# --- 1. Add the new worksheet as the last tab -----------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet  = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "GEO_RNASEQ"
$ws = $newSheet

# --- 2. Fill in the ER table content ---------------------------------------
# Column headers: <blank>, TermSourceRef, Ontology, TAN, Content type
# (validation), Notes during templating, Target term, Instruction,
# Requirement (m/o/n), Value (cv/s/d), Additional information, Review
# comments - followed by one row per column of the main sheet, giving its
# name plus (where applicable) the term source ref / ontology / term
# accession number used for that column.
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L")
$rows = @(
    @($null, 'TermSourceRef', 'Ontology', 'TAN', 'Content type (validation)', 'Notes during templating', 'Target term', 'Instruction', 'Requirement (m/o/n)', 'Value (cv/s/d)', 'Additional information', 'Review comments'),
    @('Source Name', $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null),
    @('Sample Name', $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null),
    @('Parameter [Data filtering software]', 'NFDI4PSO:0000023', 'NFDI4PSO', 'http://purl.obolibrary.org/obo/NFDI4PSO_0000023', $null, $null, $null, $null, $null, $null, $null, $null),
    @('Parameter [Data filtering software version]', 'NFDI4PSO:0000024', 'NFDI4PSO', 'http://purl.obolibrary.org/obo/NFDI4PSO_0000024', $null, $null, $null, $null, $null, $null, $null, $null),
    @('Parameter [Data filtering Software Parameters]', 'NFDI4PSO:0000025', 'NFDI4PSO', 'http://purl.obolibrary.org/obo/NFDI4PSO_0000025', $null, $null, $null, $null, $null, $null, $null, $null),
    @('Parameter [Read Alignment Software]', 'NFDI4PSO:0000002', 'NFDI4PSO', 'http://purl.obolibrary.org/obo/NFDI4PSO_0000002', $null, $null, $null, $null, $null, $null, $null, $null),
    @('Parameter [Read Alignment Software Version]', 'NFDI4PSO:0000003', 'NFDI4PSO', 'http://purl.obolibrary.org/obo/NFDI4PSO_0000003', $null, $null, $null, $null, $null, $null, $null, $null),
    @('Parameter [Read Alignment Software Parameters]', 'NFDI4PSO:0000004', 'NFDI4PSO', 'http://purl.obolibrary.org/obo/NFDI4PSO_0000004', $null, $null, $null, $null, $null, $null, $null, $null),
    @('Parameter [Genome reference sequence]', 'NFDI4PSO:0000026', 'NFDI4PSO', 'http://purl.obolibrary.org/obo/NFDI4PSO_0000026', $null, $null, $null, $null, $null, $null, $null, $null),
    @('Parameter [Processed data file name]', 'NFDI4PSO:0000028', 'NFDI4PSO', 'http://purl.obolibrary.org/obo/NFDI4PSO_0000028', $null, $null, $null, $null, $null, $null, $null, $null),
    @('Parameter [Processed data file format]', 'NFDI4PSO:0000027', 'NFDI4PSO', 'http://purl.obolibrary.org/obo/NFDI4PSO_0000027', $null, $null, $null, $null, $null, $null, $null, $null),
    @('Data File Name', $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
)

for ($r = 0; $r -lt $rows.Count; $r++) {
    $rowData = $rows[$r]
    for ($c = 0; $c -lt $cols.Count; $c++) {
        $val = $rowData[$c]
        if ($val -ne $null) {
            $ws.Range($cols[$c] + ($r + 1)).Value = $val
        }
    }
}

# --- 3. Size the columns to fit their content (bestFit widths) -------------
$widths = @(45.333333333333336, 16.666666666666668, 9.0, 45.666666666666664, 22.666666666666668, 22.0, 10.5, 9.666666666666666, 19.666666666666668, 12.666666666666666, 20.666666666666668, 16.666666666666668)
for ($c = 0; $c -lt $cols.Count; $c++) {
    $ws.Columns.Item($c + 1).ColumnWidth = $widths[$c]
}

# --- 4. Select the whole sheet (mirrors the saved selection state) and make
#        it the active tab ---------------------------------------------------
$ws.Cells.Select()
$ws.Activate()
